# Apply the "output generated at 456a3b4" data refresh to 杭州-漫展信息.xlsx
#
# Sheet map (by tab order):
#   1 = 展览      (Exhibitions)
#   2 = 演出      (Performances)
#   3 = 本地生活  (Local life)
#   4 = 全部类型  (All types - merged view of the other three sheets)

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

# ---------------------------------------------------------------------------
# Sheet1 (展览): bump "interested" counters in column F
# ---------------------------------------------------------------------------
$ws1.Range("F3").Value2  = 2600
$ws1.Range("F7").Value2  = 1903
$ws1.Range("F8").Value2  = 1741
$ws1.Range("F11").Value2 = 2411
$ws1.Range("F12").Value2 = 513
$ws1.Range("F13").Value2 = 192
$ws1.Range("F18").Value2 = 8816
$ws1.Range("F20").Value2 = 6871
$ws1.Range("F21").Value2 = 11159
$ws1.Range("F25").Value2 = 307
$ws1.Range("F26").Value2 = 532
$ws1.Range("F27").Value2 = 2453
$ws1.Range("F28").Value2 = 207
$ws1.Range("F30").Value2 = 2313
$ws1.Range("F31").Value2 = 470
$ws1.Range("F33").Value2 = 4466
$ws1.Range("F34").Value2 = 670
$ws1.Range("F35").Value2 = 300
$ws1.Range("F36").Value2 = 24
$ws1.Range("F37").Value2 = 467

# Row 23 (杭州·巨人only同人展中学篇) got cancelled: title annotated, price
# marked un-sellable (turns a numeric price into a text marker).
$ws1.Range("C23").Value2 = "杭州·巨人only同人展中学篇（取消）"
$ws1.Range("G23").Value2 = "不可售"

# ---------------------------------------------------------------------------
# Sheet2 (演出)
# ---------------------------------------------------------------------------
$ws2.Range("F8").Value2 = 1177

# ---------------------------------------------------------------------------
# Sheet3 (本地生活)
# ---------------------------------------------------------------------------
$ws3.Range("F3").Value2 = 617
$ws3.Range("F5").Value2 = 105

# ---------------------------------------------------------------------------
# Sheet4 (全部类型): same counter bumps as above (this sheet mirrors rows
# from the other three sheets).
# ---------------------------------------------------------------------------
$ws4.Range("F3").Value2  = 617
$ws4.Range("F5").Value2  = 105
$ws4.Range("F7").Value2  = 2600
$ws4.Range("F11").Value2 = 1903
$ws4.Range("F13").Value2 = 1741
$ws4.Range("F16").Value2 = 2411
$ws4.Range("F18").Value2 = 513
$ws4.Range("F19").Value2 = 192
$ws4.Range("F24").Value2 = 8816
$ws4.Range("F26").Value2 = 6871
$ws4.Range("F27").Value2 = 11159
$ws4.Range("F37").Value2 = 207
$ws4.Range("F40").Value2 = 4466
$ws4.Range("F46").Value2 = 467

# Rows 30-34 on sheet4 ripple forward by one event each (the cancelled
# "巨人only" show drops out of this merged view) and a brand-new concert
# ("法国姐姐" Joyce Jonathan tour) is appended at row 34.
#
# New row 30 (was old row 31 content)
$ws4.Range("B30").Value2 = "'2024-11-10"
$ws4.Range("C30").Value2 = "杭州·HD02动漫展嘉宾内场——锦鲤"
$ws4.Range("D30").Value2 = "钱江世纪城奔竞大道353号 杭州国际博览中心"
$ws4.Range("E30").Value2 = "2024.11.10 09:30-11.10 17:00"
$ws4.Range("F30").Value2 = 222
$ws4.Range("G30").Value2 = 258
$ws4.Range("H30").Value2 = "https://show.bilibili.com/platform/detail.html?id=92734"
$ws4.Range("I30").Value2 = "//i0.hdslb.com/bfs/openplatform/202409/NZiRZbKN1727164629427.png"

# New row 31 (was old row 32 content) - date (B31) unchanged
$ws4.Range("C31").Value2 = "杭州·崩坏同人ONLY 爱莉希雅生日会"
$ws4.Range("D31").Value2 = "康候圣街99号 顺丰创新中心"
$ws4.Range("E31").Value2 = "2024.11.10 08:00-11.10 20:00"
$ws4.Range("F31").Value2 = 307
$ws4.Range("G31").Value2 = 79
$ws4.Range("H31").Value2 = "https://show.bilibili.com/platform/detail.html?id=92228"
$ws4.Range("I31").Value2 = "//i0.hdslb.com/bfs/openplatform/202409/1FsO31h71725897488610.jpeg"

# New row 32 (was old row 33 content)
$ws4.Range("B32").Value2 = "'2024-11-16"
$ws4.Range("C32").Value2 = "杭州·1+1≥2 X PianoLab《琴键间的诗和远方》 丁阳钢琴独奏音乐会"
$ws4.Range("D32").Value2 = "江干区新业路39号 杭州大剧院"
$ws4.Range("E32").Value2 = "2024.11.16 19:30-11.16 21:00"
$ws4.Range("F32").Value2 = 0
$ws4.Range("G32").Value2 = 80
$ws4.Range("H32").Value2 = "https://show.bilibili.com/platform/detail.html?id=92325"
$ws4.Range("I32").Value2 = "//i1.hdslb.com/bfs/openplatform/202409/LrjZE4er1726039894899.jpeg"

# New row 33 (was old row 34 content) - date (B33) unchanged
$ws4.Range("C33").Value2 = "杭州·ET金色齿轮国乙同人only"
$ws4.Range("D33").Value2 = "转塘街道珊瑚沙东路9号 杭州白金汉爵大酒店"
$ws4.Range("E33").Value2 = "2024.11.16 09:30-11.16 22:00"
$ws4.Range("F33").Value2 = 532
$ws4.Range("G33").Value2 = 25
$ws4.Range("H33").Value2 = "https://show.bilibili.com/platform/detail.html?id=92511"
$ws4.Range("I33").Value2 = "//i1.hdslb.com/bfs/openplatform/202409/XfT00A611726134427042.jpeg"

# New row 34 - brand-new entry, date (B34) unchanged ("2024-11-16")
$ws4.Range("C34").Value2 = "杭州·“法国姐姐”乔伊丝·乔纳森《小意思》2024巡回演唱会【特邀嘉宾陈丽君】"
$ws4.Range("D34").Value2 = "杭州市西湖区省府路9号 浙江省人民大会堂"
$ws4.Range("E34").Value2 = "2024.11.16 19:30-11.16 21:30"
$ws4.Range("F34").Value2 = 5
$ws4.Range("G34").Value2 = 280
$ws4.Range("H34").Value2 = "https://show.bilibili.com/platform/detail.html?id=92078"
$ws4.Range("I34").Value2 = "//i2.hdslb.com/bfs/openplatform/202409/AE6VYTdf1725614295764.jpeg"
